$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-09-24 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-09-25 Wednesday", 2)

# Update the division-problem table cells in place, by (row, column) position,
# so duplicate cell text elsewhere in the table isn't touched incorrectly.
$t = $d.Tables(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "24÷3=8, 0" },
    @{ Row = 1;  Col = 2; Text = "92÷2=46, 0" },
    @{ Row = 1;  Col = 3; Text = "23÷7=3, 2" },
    @{ Row = 1;  Col = 4; Text = "48÷7=6, 6" },
    @{ Row = 1;  Col = 5; Text = "94÷9=10, 4" },

    @{ Row = 5;  Col = 1; Text = "30÷7=4, 2" },
    @{ Row = 5;  Col = 2; Text = "75÷2=37, 1" },
    @{ Row = 5;  Col = 3; Text = "39÷8=4, 7" },
    @{ Row = 5;  Col = 4; Text = "36÷9=4, 0" },
    @{ Row = 5;  Col = 5; Text = "78÷9=8, 6" },

    @{ Row = 9;  Col = 1; Text = "72÷3=24, 0" },
    @{ Row = 9;  Col = 2; Text = "72÷9=8, 0" },
    @{ Row = 9;  Col = 3; Text = "24÷5=4, 4" },
    @{ Row = 9;  Col = 4; Text = "64÷5=12, 4" },
    @{ Row = 9;  Col = 5; Text = "21÷7=3, 0" },

    @{ Row = 13; Col = 1; Text = "52÷5=10, 2" },
    @{ Row = 13; Col = 2; Text = "62÷6=10, 2" },
    @{ Row = 13; Col = 3; Text = "15÷6=2, 3" },
    @{ Row = 13; Col = 4; Text = "87÷4=21, 3" },
    @{ Row = 13; Col = 5; Text = "51÷8=6, 3" },

    @{ Row = 17; Col = 1; Text = "79÷2=39, 1" },
    @{ Row = 17; Col = 2; Text = "55÷3=18, 1" },
    @{ Row = 17; Col = 3; Text = "73÷3=24, 1" },
    @{ Row = 17; Col = 4; Text = "41÷4=10, 1" },
    @{ Row = 17; Col = 5; Text = "93÷5=18, 3" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
